$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number (would otherwise be
# auto-converted to a numeric value by Excel, losing formatting like trailing
# zeros). Force them to remain text, then restore the default "Normal" style
# so no visible formatting change is left behind.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.990"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.259"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.566"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.992"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.558"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0504"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "53.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "88.29"
$ws.Range("D50").Style = "Normal"

# Remaining cells are unambiguous text (percentages, or numbers with multiple
# "." separators) and can simply be assigned.
$ws.Range("D2").Value = "29.852.85"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  -0.97%  "
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("E8").Value = "  +8.71%  "
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.850.67"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "1.635.77"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("E14").Value = "  +5.59%  "
$ws.Range("E15").Value = "  +5.05%  "
$ws.Range("D16").Value = "29.832.52"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("E17").Value = "  +16.94%  "
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("E22").Value = "  +2.44%  "
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +2.39%  "
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("E28").Value = "  +3.19%  "
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("E31").Value = "  +4.47%  "
$ws.Range("E32").Value = "  +3.25%  "
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").Value = "1.414.86"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  +6.84%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("E40").Value = "  +3.72%  "
$ws.Range("E41").Value = "  +3.92%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("E45").Value = "  +5.64%  "
$ws.Range("E46").Value = "  +18.22%  "
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("E48").Value = "  +2.51%  "
$ws.Range("D49").Value = "1.759.64"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("E51").Value = "  +2.07%  "
